$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: add trailing spaces, then append a red parenthetical
#    note split across three runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs.First
$insertPos = $p1.Range.End - 1

$t1 = "(This is a change " + [char]0x2013 + " Ve"
$t2 = "rsion for main branch"
$t3 = ")"

$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter($t1)
$len1 = $t1.Length
$c1 = $d.Range($insertPos, $insertPos + $len1)
$c1.Font.Color = 255

$pos2 = $insertPos + $len1
$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter($t2)
$len2 = $t2.Length
$c2 = $d.Range($pos2, $pos2 + $len2)
$c2.Font.Color = 255

$pos3 = $pos2 + $len2
$r3 = $d.Range($pos3, $pos3)
$r3.InsertAfter($t3)
$len3 = $t3.Length
$c3 = $d.Range($pos3, $pos3 + $len3)
$c3.Font.Color = 255

# ---------------------------------------------------------------------------
# 2) "Crispian's Day speech" paragraph: merge the floating space into the
#    preceding run, and merge the " Henry V" .. "]" runs into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Day speech from ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Day speech from ", 2) | Out-Null

$d.Content.Find.Execute("Henry V [Source " + [char]0x2013 + " Wikipedia]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Henry V [Source " + [char]0x2013 + " Wikipedia]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append a new, empty "larger" styled paragraph after the final
#    "... Saint Crispin's day." paragraph.
# ---------------------------------------------------------------------------
$tail = $d.Content
$tail.Find.Execute("Saint Crispin(')s day.", $true, $false, $true, $false, $false,
                    $true, 1, $false, "Saint Crispin\1s day.^p", 2) | Out-Null

$newPara = $d.Paragraphs.Last
$newPara.Style = "larger"
$newPara.Format.Shading.Texture = 0
$newPara.Format.Shading.ForegroundPatternColor = -16777216
$newPara.Format.Shading.BackgroundPatternColor = 16777215
$newPara.Format.SpaceBefore = 0
$newPara.Format.SpaceBeforeAuto = $false
$newPara.Format.SpaceAfter = 7.5
$newPara.Format.SpaceAfterAuto = $false

# ---------------------------------------------------------------------------
# 4) Styles cleanup: drop the unused "apple-converted-space" and
#    "Hyperlink" character styles.
# ---------------------------------------------------------------------------
for ($i = $d.Styles.Count; $i -ge 1; $i--) {
    $nm = $d.Styles.Item($i).NameLocal
    if ($nm -eq "apple-converted-space" -or $nm -eq "Hyperlink") {
        $d.Styles.Item($i).Delete()
    }
}

Write-Output "done"
